$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D/E/F (and some C) values for rows 2-23 per the corrected dataset.
$data = @{
    2  = @{ C = 1252; D = 861;  E = 937;  F = 226 }
    3  = @{ C = 1187; D = 699;  E = 855;  F = 210 }
    4  = @{ C = 1324; D = 863;  E = 687;  F = 255 }
    5  = @{ C = 1517; D = 726;  E = 693;  F = 255 }
    6  = @{ C = 925;  D = 837;  E = 736;  F = 365 }
    7  = @{ C = 404;  D = 531;  E = 531;  F = 286 }
    8  = @{ C = 584;  D = 728;  E = 855;  F = 333 }
    9  = @{ C = 686;  D = 709;  E = 875;  F = 390 }
    10 = @{ C = 641;  D = 731;  E = 857;  F = 450 }
    11 = @{ C = 832;  D = 990;  E = 1202; F = 567 }
    12 = @{ C = 868;  D = 1196; E = 1101; F = 452 }
    13 = @{ C = 1197; D = 1561; E = 1373; F = 594 }
    14 = @{ C = 1089; D = 1651; E = 1364; F = 582 }
    15 = @{ C = 1225; D = 1891; E = 1512; F = 627 }
    16 = @{ C = 1232; D = 1918; E = 1395; F = 573 }
    17 = @{ C = 1277; D = 2013; E = 1415; F = 512 }
    18 = @{ C = 1125; D = 1573; E = 1268; F = 532 }
    19 = @{ C = 1043; D = 1579; E = 1188; F = 520 }
    20 = @{ C = 1159; D = 1845; E = 1368; F = 594 }
    21 = @{ C = 1288; D = 1826; E = 1407; F = 520 }
    22 = @{ C = 1225; D = 1858; E = 1329; F = 581 }
    23 = @{ C = 1208; D = 2096; E = 1330; F = 630 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
}

# Update the active selection to match the saved workbook view (F24).
$ws.Range("F24").Select()
